# Apply updated '想去人数' (want-to-go count) values per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1819
$ws.Range("F8").Value = 629
$ws.Range("F9").Value = 3608
$ws.Range("F10").Value = 1219
$ws.Range("F15").Value = 1412
$ws.Range("F17").Value = 1803
$ws.Range("F20").Value = 7
$ws.Range("F27").Value = 276
$ws.Range("F29").Value = 4379
$ws.Range("F31").Value = 11
$ws.Range("F35").Value = 1225
$ws.Range("F37").Value = 11

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F21").Value = 94
$ws.Range("F23").Value = 150
$ws.Range("F36").Value = 452
$ws.Range("F43").Value = 91

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F9").Value = 397
$ws.Range("F10").Value = 3026
$ws.Range("F11").Value = 534
$ws.Range("F12").Value = 825
$ws.Range("F13").Value = 275
$ws.Range("F14").Value = 277

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 534
$ws.Range("F8").Value = 825
$ws.Range("F9").Value = 275
$ws.Range("F12").Value = 629
$ws.Range("F13").Value = 3608
$ws.Range("F14").Value = 1219
$ws.Range("F18").Value = 1412
$ws.Range("F23").Value = 1803
$ws.Range("F25").Value = 7
$ws.Range("F27").Value = 94
$ws.Range("F30").Value = 150
$ws.Range("F31").Value = 150
$ws.Range("F37").Value = 276
$ws.Range("F39").Value = 277
$ws.Range("F41").Value = 4379
$ws.Range("F42").Value = 452
$ws.Range("F43").Value = 11
$ws.Range("F46").Value = 91
$ws.Range("F50").Value = 1225
$ws.Range("F52").Value = 11
